$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.373.72"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3
$ws.Range("D3").Value = "1.898.33"
$ws.Range("E3").Value = "  -1.14%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.34%  "

# Row 5
$ws.Range("D5").Value = "'325.07"
$ws.Range("E5").Value = "  -3.00%  "

# Row 6
$ws.Range("E6").Value = "  +0.15%  "

# Row 7
$ws.Range("D7").Value = "'0.4796"
$ws.Range("E7").Value = "  +2.51%  "

# Row 8
$ws.Range("D8").Value = "'0.4060"
$ws.Range("E8").Value = "  -1.24%  "

# Row 9
$ws.Range("D9").Value = "'0.08063"
$ws.Range("E9").Value = "  +0.32%  "

# Row 10
$ws.Range("E10").Value = "  -1.63%  "

# Row 11
$ws.Range("D11").Value = "'23.22"
$ws.Range("E11").Value = "  +3.42%  "

# Row 12
$ws.Range("D12").Value = "1.996.84"
$ws.Range("E12").Value = "  +3.56%  "

# Row 13
$ws.Range("D13").Value = "'5.943"
$ws.Range("E13").Value = "  -0.93%  "

# Row 14
$ws.Range("E14").Value = "  -1.74%  "

# Row 15
$ws.Range("D15").Value = "'89.81"
$ws.Range("E15").Value = "  -0.25%  "

# Row 16
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.24%  "

# Row 17
$ws.Range("D17").Value = "'0.06683"
$ws.Range("E17").Value = "  +1.37%  "

# Row 18
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
$ws.Range("D19").Value = "'17.63"
$ws.Range("E19").Value = "  -1.28%  "

# Row 20
$ws.Range("E20").Value = "  +0.17%  "

# Row 21
$ws.Range("D21").Value = "29.387.06"
$ws.Range("E21").Value = "  -0.87%  "

# Row 22
$ws.Range("D22").Value = "'5.530"
$ws.Range("E22").Value = "  -0.95%  "

# Row 23
$ws.Range("D23").Value = "'11.72"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("D24").Value = "'2.160"
$ws.Range("E24").Value = "  -2.24%  "

# Row 25
$ws.Range("D25").Value = "2.144.10"
$ws.Range("E25").Value = "  -1.24%  "

# Row 26
$ws.Range("D26").Value = "'155.02"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27
$ws.Range("D27").Value = "'19.75"
$ws.Range("E27").Value = "  -0.63%  "

# Row 28
$ws.Range("D28").Value = "'6.095"
$ws.Range("E28").Value = "  +6.31%  "

# Row 29
$ws.Range("D29").Value = "'2.089"
$ws.Range("E29").Value = "  -2.66%  "

# Row 30
$ws.Range("D30").Value = "'118.19"
$ws.Range("E30").Value = "  +0.54%  "

# Row 31
$ws.Range("D31").Value = "'1.021"
$ws.Range("E31").Value = "  -4.47%  "

# Row 32
$ws.Range("D32").Value = "'0.09514"
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$ws.Range("D33").Value = "'1.388"
$ws.Range("E33").Value = "  -3.43%  "

# Row 34
$ws.Range("D34").Value = "'3.527"
$ws.Range("E34").Value = "  -1.32%  "

# Row 35
$ws.Range("D35").Value = "'5.385"
$ws.Range("E35").Value = "  -0.72%  "

# Row 36
$ws.Range("E36").Value = "  -1.09%  "

# Row 37
$ws.Range("D37").Value = "'0.06056"
$ws.Range("E37").Value = "  -1.43%  "

# Row 38
$ws.Range("D38").Value = "'1.174"
$ws.Range("E38").Value = "  -0.75%  "

# Row 39
$ws.Range("D39").Value = "'0.5861"
$ws.Range("E39").Value = "  -0.54%  "

# Row 40
$ws.Range("D40").Value = "'7.864"
$ws.Range("E40").Value = "  -6.89%  "

# Row 41
$ws.Range("D41").Value = "'0.1844"
$ws.Range("E41").Value = "  -0.17%  "

# Row 42
$ws.Range("D42").Value = "'10.18"
$ws.Range("E42").Value = "  -0.57%  "

# Row 43
$ws.Range("D43").Value = "'1.288"
$ws.Range("E43").Value = "  +2.12%  "

# Row 44
$ws.Range("D44").Value = "'2.403"
$ws.Range("E44").Value = "  +1.91%  "

# Row 45
$ws.Range("D45").Value = "'0.07718"
$ws.Range("E45").Value = "  +2.76%  "

# Row 46
$ws.Range("D46").Value = "'12.19"
$ws.Range("E46").Value = "  -0.91%  "

# Row 47
$ws.Range("D47").Value = "'0.5512"
$ws.Range("E47").Value = "  -1.36%  "

# Row 48
$ws.Range("D48").Value = "'1.923"
$ws.Range("E48").Value = "  -0.69%  "

# Row 49
$ws.Range("D49").Value = "'113.06"
$ws.Range("E49").Value = "  -0.21%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'72.11"
$ws.Range("E50").Value = "  +0.85%  "

# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.2926"
$ws.Range("E51").Value = "  -2.48%  "
